# Auto-generated Excel COM-interop script applying the Typhon_Profits
# profit-recalculation update across all 8 job sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Columns H-N hold cached market/profit figures
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) that were refreshed by the scheduled pricing runner.

$wb = $excel.ActiveWorkbook

# ===================== ALC =====================
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 198.5
$ws.Range("I33").Value = 204.35294
$ws.Range("J33").Value = 99
$ws.Range("K33").Value = 204.35294
$ws.Range("L33").Value = 99
$ws.Range("M33").Value = 24.64706000000001
$ws.Range("N33").Value = -557
# Row 98
$ws.Range("H98").Value = 1262.375
$ws.Range("I98").Value = 1080
$ws.Range("J98").Value = 1566.3334
$ws.Range("K98").Value = 1080
$ws.Range("L98").Value = 1566.3334
$ws.Range("M98").Value = 418
$ws.Range("N98").Value = -4562.3334
# Row 112
$ws.Range("H112").Value = 1096.1945
$ws.Range("J112").Value = 1113.6177
$ws.Range("L112").Value = 3340.8531
$ws.Range("N112").Value = -5556.8531
# Row 122
$ws.Range("H122").Value = 1262.375
$ws.Range("I122").Value = 1080
$ws.Range("J122").Value = 1566.3334
$ws.Range("K122").Value = 3240
$ws.Range("L122").Value = 4699.0002
$ws.Range("M122").Value = -790
$ws.Range("N122").Value = -9599.0002
# Row 129
$ws.Range("H129").Value = 1017
$ws.Range("J129").Value = 1038.8077
$ws.Range("L129").Value = 3116.4231
$ws.Range("N129").Value = -13116.4231
# Row 132
$ws.Range("H132").Value = 2721.639
$ws.Range("I132").Value = 2808.0293
$ws.Range("K132").Value = 8424.0879
$ws.Range("M132").Value = -5894.0879
# Row 141
$ws.Range("H141").Value = 2551.0588
$ws.Range("J141").Value = 2953.1667
$ws.Range("L141").Value = 8859.500100000001
$ws.Range("N141").Value = -19219.5001

# ===================== ARM =====================
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 50000000
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# Row 32
$ws.Range("H32").Value = 5094.2104
$ws.Range("I32").Value = 3880.8354
$ws.Range("J32").Value = 11085.25
$ws.Range("K32").Value = 3880.8354
$ws.Range("L32").Value = 11085.25
$ws.Range("M32").Value = -3593.8354
$ws.Range("N32").Value = -11659.25
# Row 61
$ws.Range("H61").Value = 3452.8518
$ws.Range("I61").Value = 3613.15
$ws.Range("J61").Value = 2994.8572
$ws.Range("K61").Value = 3613.15
$ws.Range("L61").Value = 2994.8572
$ws.Range("M61").Value = -3401.15
$ws.Range("N61").Value = -3418.8572
# Row 74
$ws.Range("H74").Value = 33335014
$ws.Range("I74").Value = 43478904
$ws.Range("K74").Value = 43478904
$ws.Range("M74").Value = -43478030
# Row 77
$ws.Range("H77").Value = 33335014
$ws.Range("I77").Value = 43478904
$ws.Range("K77").Value = 217394520
$ws.Range("M77").Value = -217390152
# Row 134
$ws.Range("H134").Value = 59500
$ws.Range("J134").Value = 59500
$ws.Range("L134").Value = 59500
$ws.Range("N134").Value = -69640
# Row 136
$ws.Range("H136").Value = 3452.8518
$ws.Range("I136").Value = 3613.15
$ws.Range("J136").Value = 2994.8572
$ws.Range("K136").Value = 10839.45
$ws.Range("L136").Value = 8984.571599999999
$ws.Range("M136").Value = -8289.450000000001
$ws.Range("N136").Value = -14084.5716

# ===================== BSM =====================
$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 4333
$ws.Range("I11").Value = 4333
$ws.Range("K11").Value = 4333
$ws.Range("M11").Value = -4193
# Row 20
$ws.Range("H20").Value = 2118.6365
$ws.Range("I20").Value = 2353.6428
$ws.Range("J20").Value = 1707.375
$ws.Range("K20").Value = 2353.6428
$ws.Range("L20").Value = 1707.375
$ws.Range("M20").Value = -2106.6428
$ws.Range("N20").Value = -2201.375
# Row 105
$ws.Range("H105").Value = 1825.9584
$ws.Range("I105").Value = 1611.2307
$ws.Range("K105").Value = 1611.2307
$ws.Range("M105").Value = 135.7692999999999
# Row 134
$ws.Range("H134").Value = 3217.0488
$ws.Range("I134").Value = 3255.9119
$ws.Range("K134").Value = 9767.735700000001
$ws.Range("M134").Value = -7232.735700000001

# ===================== CRP =====================
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 31
$ws.Range("H31").Value = 3609.6487
$ws.Range("I31").Value = 2788
$ws.Range("K31").Value = 2788
$ws.Range("M31").Value = -2493
# Row 34
$ws.Range("H34").Value = 3609.6487
$ws.Range("I34").Value = 2788
$ws.Range("K34").Value = 2788
$ws.Range("M34").Value = -2586
# Row 58
$ws.Range("H58").Value = 24563.137
$ws.Range("I58").Value = 1878.909
$ws.Range("J58").Value = 47247.363
$ws.Range("K58").Value = 1878.909
$ws.Range("L58").Value = 47247.363
$ws.Range("M58").Value = -1675.909
$ws.Range("N58").Value = -47653.363
# Row 122
$ws.Range("H122").Value = 1625
$ws.Range("I122").Value = 1550
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4650
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2200
$ws.Range("N122").Value = -10900
# Row 134
$ws.Range("H134").Value = 1302.2727
$ws.Range("I134").Value = 1115.5
$ws.Range("K134").Value = 3346.5
$ws.Range("M134").Value = -811.5
# Row 136
$ws.Range("H136").Value = 24563.137
$ws.Range("I136").Value = 1878.909
$ws.Range("J136").Value = 47247.363
$ws.Range("K136").Value = 5636.727000000001
$ws.Range("L136").Value = 141742.089
$ws.Range("M136").Value = -3086.727000000001
$ws.Range("N136").Value = -146842.089

# ===================== CUL =====================
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 62500130
$ws.Range("I38").Value = 70
$ws.Range("J38").Value = 100000160
$ws.Range("K38").Value = 210
$ws.Range("L38").Value = 300000480
$ws.Range("M38").Value = 137
$ws.Range("N38").Value = -300001174
# Row 99
$ws.Range("H99").Value = 3610.15
$ws.Range("I99").Value = 1923.3846
$ws.Range("J99").Value = 6742.7144
$ws.Range("K99").Value = 5770.1538
$ws.Range("L99").Value = 20228.1432
$ws.Range("M99").Value = -3524.1538
$ws.Range("N99").Value = -24720.1432
# Row 121
$ws.Range("H121").Value = 1055.84
$ws.Range("J121").Value = 1087.3334
$ws.Range("L121").Value = 3262.0002
$ws.Range("N121").Value = -5882.0002
# Row 131
$ws.Range("H131").Value = 701.59
$ws.Range("J131").Value = 701.59
$ws.Range("L131").Value = 2104.77
$ws.Range("N131").Value = -12184.77
# Row 137
$ws.Range("H137").Value = 17548744
$ws.Range("J137").Value = 20838796
$ws.Range("L137").Value = 62516388
$ws.Range("N137").Value = -62526588

# ===================== GSM =====================
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 14168.6
$ws.Range("I70").Value = 6500
$ws.Range("K70").Value = 6500
$ws.Range("M70").Value = -6230
# Row 73
$ws.Range("H73").Value = 14168.6
$ws.Range("I73").Value = 6500
$ws.Range("K73").Value = 6500
$ws.Range("M73").Value = -5564
# Row 132
$ws.Range("H132").Value = 14530.286
$ws.Range("I132").Value = 2648
$ws.Range("J132").Value = 169000
$ws.Range("K132").Value = 7944
$ws.Range("L132").Value = 507000
$ws.Range("M132").Value = -5414
$ws.Range("N132").Value = -512060
# Row 134
$ws.Range("H134").Value = 22000
$ws.Range("J134").Value = 22000
$ws.Range("L134").Value = 66000
$ws.Range("N134").Value = -71070

# ===================== LTW =====================
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2999
$ws.Range("J68").Value = 2999
$ws.Range("L68").Value = 2999
$ws.Range("N68").Value = -4497
# Row 71
$ws.Range("H71").Value = 2999
$ws.Range("J71").Value = 2999
$ws.Range("L71").Value = 14995
$ws.Range("N71").Value = -22483
# Row 135
$ws.Range("H135").Value = 30214.5
$ws.Range("J135").Value = 30214.5
$ws.Range("L135").Value = 30214.5
$ws.Range("N135").Value = -40354.5
# Row 136
$ws.Range("H136").Value = 1385.4062
$ws.Range("I136").Value = 1288.2307
$ws.Range("J136").Value = 1806.5
$ws.Range("K136").Value = 3864.6921
$ws.Range("L136").Value = 5419.5
$ws.Range("M136").Value = -1314.6921
$ws.Range("N136").Value = -10519.5

# ===================== WVR =====================
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1378.7273
$ws.Range("I132").Value = 1175.3684
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 3526.1052
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -996.1052
$ws.Range("N132").Value = -13060.0001
